$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 204; everything currently at 204:257 shifts down to 205:258,
# and Excel auto-extends the used range / dimension to R258.
$ws.Rows("204:204").Insert()

# Populate the newly inserted row 204 with the new weekly record.
$ws.Cells.Item(204, 1).Value = 3
$ws.Cells.Item(204, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(204, 3).Value = "Coquimbo"
$ws.Cells.Item(204, 4).Value = 44551
$ws.Cells.Item(204, 5).Value = 5
$ws.Cells.Item(204, 6).Value = 100112043
$ws.Cells.Item(204, 7).Value = "Pepino ensalada"
$ws.Cells.Item(204, 8).Value = "Sin especificar"
$ws.Cells.Item(204, 9).Value = "Primera"
$ws.Cells.Item(204, 10).Value = 105
$ws.Cells.Item(204, 11).Value = 7500
$ws.Cells.Item(204, 12).Value = 8000
$ws.Cells.Item(204, 13).Value = 7762
$ws.Cells.Item(204, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(204, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(204, 16).Value = 111
$ws.Cells.Item(204, 17).Value = 70
$ws.Cells.Item(204, 18).Value = "Hortaliza"
